$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.180165333333334
$ws.Range("H2").Value = 6.540496
$ws.Range("I2").Value = 0.01970539991828544
$ws.Range("J2").Value = 0.01970539991828544
$ws.Range("M2").Value = 1.815761
$ws.Range("N2").Value = 5.447283000000001
$ws.Range("O2").Value = 0.07007596730428067
$ws.Range("P2").Value = 0.07007596730428067
$ws.Range("Q2").Value = 3.958659185818667
$ws.Range("R2").Value = 35.627932672368
$ws.Range("S2").Value = 0.001380874960391545
$ws.Range("T2").Value = 0.001380874960391545

$ws.Range("G3").Value = 2.180165333333334
$ws.Range("H3").Value = 6.540496
$ws.Range("I3").Value = 0.01970539991828544
$ws.Range("J3").Value = 0.01970539991828544
$ws.Range("O3").Value = 0.5079540516959071
$ws.Range("P3").Value = 0.5079540516959072
$ws.Range("Q3").Value = 28.69481578454045
$ws.Range("R3").Value = 258.253342060864
$ws.Range("S3").Value = 0.01000943772878128
$ws.Range("T3").Value = 0.01000943772878129

$ws.Range("G4").Value = 2.180165333333334
$ws.Range("H4").Value = 6.540496
$ws.Range("I4").Value = 0.01970539991828544
$ws.Range("J4").Value = 0.01970539991828544
$ws.Range("M4").Value = 9.711409333333334
$ws.Range("N4").Value = 29.134228
$ws.Range("O4").Value = 0.3747940411327002
$ws.Range("P4").Value = 0.3747940411327002
$ws.Range("Q4").Value = 21.17247796634311
$ws.Range("R4").Value = 190.552301697088
$ws.Range("S4").Value = 0.007385466467510178
$ws.Range("T4").Value = 0.007385466467510178

$ws.Range("G5").Value = 2.180165333333334
$ws.Range("H5").Value = 6.540496
$ws.Range("I5").Value = 0.01970539991828544
$ws.Range("J5").Value = 0.01970539991828544
$ws.Range("M5").Value = 1.222391
$ws.Range("N5").Value = 3.667173
$ws.Range("O5").Value = 0.04717593986711188
$ws.Range("P5").Value = 0.04717593986711189
$ws.Range("Q5").Value = 2.665014481978667
$ws.Range("R5").Value = 23.985130337808
$ws.Range("S5").Value = 0.0009296207616024251
$ws.Range("T5").Value = 0.0009296207616024253

$ws.Range("I6").Value = 0.733713204346044
$ws.Range("J6").Value = 0.7337132043460441
$ws.Range("M6").Value = 1.815761
$ws.Range("N6").Value = 5.447283000000001
$ws.Range("O6").Value = 0.07007596730428067
$ws.Range("P6").Value = 0.07007596730428067
$ws.Range("Q6").Value = 147.3971869733887
$ws.Range("R6").Value = 1326.574682760498
$ws.Range("S6").Value = 0.05141566251847238
$ws.Range("T6").Value = 0.05141566251847239

$ws.Range("I7").Value = 0.733713204346044
$ws.Range("J7").Value = 0.7337132043460441
$ws.Range("O7").Value = 0.5079540516959071
$ws.Range("P7").Value = 0.5079540516959072
$ws.Range("S7").Value = 0.3726925949303601
$ws.Range("T7").Value = 0.3726925949303603

$ws.Range("I8").Value = 0.733713204346044
$ws.Range("J8").Value = 0.7337132043460441
$ws.Range("M8").Value = 9.711409333333334
$ws.Range("N8").Value = 29.134228
$ws.Range("O8").Value = 0.3747940411327002
$ws.Range("P8").Value = 0.3747940411327002
$ws.Range("Q8").Value = 788.3385628837965
$ws.Range("R8").Value = 7095.047065954168
$ws.Range("S8").Value = 0.2749913368892765
$ws.Range("T8").Value = 0.2749913368892765

$ws.Range("I9").Value = 0.733713204346044
$ws.Range("J9").Value = 0.7337132043460441
$ws.Range("M9").Value = 1.222391
$ws.Range("N9").Value = 3.667173
$ws.Range("O9").Value = 0.04717593986711188
$ws.Range("P9").Value = 0.04717593986711189
$ws.Range("Q9").Value = 99.22946620264867
$ws.Range("R9").Value = 893.0651958238381
$ws.Range("S9").Value = 0.03461361000793495
$ws.Range("T9").Value = 0.03461361000793495

$ws.Range("G10").Value = 25.672264
$ws.Range("H10").Value = 77.016792
$ws.Range("I10").Value = 0.2320384702908474
$ws.Range("J10").Value = 0.2320384702908474
$ws.Range("M10").Value = 1.815761
$ws.Range("N10").Value = 5.447283000000001
$ws.Range("O10").Value = 0.07007596730428067
$ws.Range("P10").Value = 0.07007596730428067
$ws.Range("Q10").Value = 46.61469575290401
$ws.Range("R10").Value = 419.532261776136
$ws.Range("S10").Value = 0.01626032025743672
$ws.Range("T10").Value = 0.01626032025743672

$ws.Range("G11").Value = 25.672264
$ws.Range("H11").Value = 77.016792
$ws.Range("I11").Value = 0.2320384702908474
$ws.Range("J11").Value = 0.2320384702908474
$ws.Range("O11").Value = 0.5079540516959071
$ws.Range("P11").Value = 0.5079540516959072
$ws.Range("Q11").Value = 337.8922116543253
$ws.Range("R11").Value = 3041.029904888928
$ws.Range("S11").Value = 0.1178648811335563
$ws.Range("T11").Value = 0.1178648811335563

$ws.Range("G12").Value = 25.672264
$ws.Range("H12").Value = 77.016792
$ws.Range("I12").Value = 0.2320384702908474
$ws.Range("J12").Value = 0.2320384702908474
$ws.Range("M12").Value = 9.711409333333334
$ws.Range("N12").Value = 29.134228
$ws.Range("O12").Value = 0.3747940411327002
$ws.Range("P12").Value = 0.3747940411327002
$ws.Range("Q12").Value = 249.3138642173973
$ws.Range("R12").Value = 2243.824777956576
$ws.Range("S12").Value = 0.08696663597855668
$ws.Range("T12").Value = 0.08696663597855668

$ws.Range("G13").Value = 25.672264
$ws.Range("H13").Value = 77.016792
$ws.Range("I13").Value = 0.2320384702908474
$ws.Range("J13").Value = 0.2320384702908474
$ws.Range("M13").Value = 1.222391
$ws.Range("N13").Value = 3.667173
$ws.Range("O13").Value = 0.04717593986711188
$ws.Range("P13").Value = 0.04717593986711189
$ws.Range("Q13").Value = 31.381544463224
$ws.Range("R13").Value = 282.433900169016
$ws.Range("S13").Value = 0.01094663292129764
$ws.Range("T13").Value = 0.01094663292129765

$ws.Range("G14").Value = 1.608999666666667
$ws.Range("H14").Value = 4.826999
$ws.Range("I14").Value = 0.01454292544482312
$ws.Range("J14").Value = 0.01454292544482312
$ws.Range("M14").Value = 1.815761
$ws.Range("N14").Value = 5.447283000000001
$ws.Range("O14").Value = 0.07007596730428067
$ws.Range("P14").Value = 0.07007596730428067
$ws.Range("Q14").Value = 2.921558843746333
$ws.Range("R14").Value = 26.294029593717
$ws.Range("S14").Value = 0.001019109567980017
$ws.Range("T14").Value = 0.001019109567980017

$ws.Range("G15").Value = 1.608999666666667
$ws.Range("H15").Value = 4.826999
$ws.Range("I15").Value = 0.01454292544482312
$ws.Range("J15").Value = 0.01454292544482312
$ws.Range("O15").Value = 0.5079540516959071
$ws.Range("P15").Value = 0.5079540516959072
$ws.Range("Q15").Value = 21.17726959807955
$ws.Range("R15").Value = 190.595426382716
$ws.Range("S15").Value = 0.007387137903209409
$ws.Range("T15").Value = 0.00738713790320941

$ws.Range("G16").Value = 1.608999666666667
$ws.Range("H16").Value = 4.826999
$ws.Range("I16").Value = 0.01454292544482312
$ws.Range("J16").Value = 0.01454292544482312
$ws.Range("M16").Value = 9.711409333333334
$ws.Range("N16").Value = 29.134228
$ws.Range("O16").Value = 0.3747940411327002
$ws.Range("P16").Value = 0.3747940411327002
$ws.Range("Q16").Value = 15.62565438019689
$ws.Range("R16").Value = 140.630889421772
$ws.Range("S16").Value = 0.00545060179735683
$ws.Range("T16").Value = 0.00545060179735683

$ws.Range("G17").Value = 1.608999666666667
$ws.Range("H17").Value = 4.826999
$ws.Range("I17").Value = 0.01454292544482312
$ws.Range("J17").Value = 0.01454292544482312
$ws.Range("M17").Value = 1.222391
$ws.Range("N17").Value = 3.667173
$ws.Range("O17").Value = 0.04717593986711188
$ws.Range("P17").Value = 0.04717593986711189
$ws.Range("Q17").Value = 1.966826711536333
$ws.Range("R17").Value = 17.701440403827
$ws.Range("S17").Value = 0.000686076176276867
$ws.Range("T17").Value = 0.0006860761762768671
